$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "60.841.41"
$ws.Range('E2').Value = "  -2.36%  "

$ws.Range('D3').Value = "2.419.39"
$ws.Range('E3').Value = "  -1.69%  "

$ws.Range('D4').Value = "'0.997"
$ws.Range('E4').Value = "  -0.21%  "

$ws.Range('D5').Value = "'570.90"
$ws.Range('E5').Value = "  -1.15%  "

$ws.Range('D6').Value = "'139.57"
$ws.Range('E6').Value = "  -3.42%  "

$ws.Range('E7').Value = "  +0.23%  "

$ws.Range('D8').Value = "'0.528"
$ws.Range('E8').Value = "  -1.20%  "

$ws.Range('D9').Value = "2.405.65"
$ws.Range('E9').Value = "  -2.14%  "

$ws.Range('E10').Value = "  -1.39%  "

$ws.Range('E11').Value = "  +0.10%  "

$ws.Range('D12').Value = "'5.09"
$ws.Range('E12').Value = "  -2.55%  "

$ws.Range('D13').Value = "'0.340"
$ws.Range('E13').Value = "  -1.72%  "

$ws.Range('D14').Value = "'25.96"
$ws.Range('E14').Value = "  -1.92%  "

$ws.Range('D15').Value = "'0.0000171"
$ws.Range('E15').Value = "  -2.35%  "

$ws.Range('D16').Value = "2.818.49"
$ws.Range('E16').Value = "  -2.76%  "

$ws.Range('D17').Value = "60.678.55"
$ws.Range('E17').Value = "  -2.50%  "

$ws.Range('D18').Value = "2.399.58"
$ws.Range('E18').Value = "  -2.77%  "

$ws.Range('D19').Value = "'10.61"
$ws.Range('E19').Value = "  -2.71%  "

$ws.Range('D20').Value = "'7.36"
$ws.Range('E20').Value = "  +2.65%  "

$ws.Range('D21').Value = "'322.42"
$ws.Range('E21').Value = "  -1.95%  "

$ws.Range('D22').Value = "'4.05"
$ws.Range('E22').Value = "  -2.07%  "

$ws.Range('D23').Value = "'6.06"
$ws.Range('E23').Value = "  +1.03%  "

$ws.Range('D25').Value = "'1.89"
$ws.Range('E25').Value = "  -4.33%  "

$ws.Range('D26').Value = "'64.83"
$ws.Range('E26').Value = "  -1.48%  "

$ws.Range('D27').Value = "'8.51"
$ws.Range('E27').Value = "  -8.27%  "

$ws.Range('D28').Value = "'575.23"
$ws.Range('E28').Value = "  -5.44%  "

$ws.Range('D29').Value = "2.519.26"
$ws.Range('E29').Value = "  -2.64%  "

$ws.Range('E30').Value = "  -4.93%  "

$ws.Range('D31').Value = "'7.93"
$ws.Range('E31').Value = "  -1.22%  "

$ws.Range('D32').Value = "'1.35"
$ws.Range('E32').Value = "  -6.32%  "

$ws.Range('D33').Value = "'1.84"
$ws.Range('E33').Value = "  -2.39%  "

$ws.Range('D34').Value = "'0.134"
$ws.Range('E34').Value = "  -3.55%  "

$ws.Range('E35').Value = "  +0.07%  "

$ws.Range('D36').Value = "'4.63"
$ws.Range('E36').Value = "  -6.19%  "

$ws.Range('B37').Value = "ImmutableX"
$ws.Range('C37').Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('D37').Value = "'1.40"
$ws.Range('E37').Value = "  -3.89%  "

$ws.Range('B38').Value = "PolygonEcosystemToken"
$ws.Range('C38').Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range('D38').Value = "'0.369"
$ws.Range('E38').Value = "  -2.41%  "

$ws.Range('B39').Value = "Monero"
$ws.Range('C39').Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D39').Value = "'149.69"
$ws.Range('E39').Value = "  -0.66%  "

$ws.Range('D40').Value = "'18.25"
$ws.Range('E40').Value = "  -1.24%  "

$ws.Range('D41').Value = "'5.14"
$ws.Range('E41').Value = "  -4.67%  "

$ws.Range('D43').Value = "'1.67"
$ws.Range('E43').Value = "  -3.79%  "

$ws.Range('D44').Value = "'41.03"
$ws.Range('E44').Value = "  -4.62%  "

$ws.Range('E45').Value = "  -7.86%  "

$ws.Range('D46').Value = "0.0₆0277"
$ws.Range('E46').Value = "  +8.69%  "

$ws.Range('D47').Value = "'141.12"
$ws.Range('E47').Value = "  -1.62%  "

$ws.Range('D48').Value = "'3.52"
$ws.Range('E48').Value = "  -3.33%  "

$ws.Range('D49').Value = "'0.589"
$ws.Range('E49').Value = "  -2.99%  "

$ws.Range('D50').Value = "'0.0504"
$ws.Range('E50').Value = "  -4.36%  "

$ws.Range('D51').Value = "'19.41"
$ws.Range('E51').Value = "  -1.89%  "
